$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Carl Montpetit ()" -> "Carl Montpetit (MONC08069000)"
#    split into three runs: "Carl Montpetit (" / "MONC08069000" / ")"
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.Execute("Carl Montpetit (")
$rng1.Collapse(0)
$rng1.InsertAfter("MONC08069000")
# force a run split from the identically-formatted neighbours, then remove
# the temporary marker so the run properties end up identical again
$rng1.Font.Bold = 1
$rng1.Font.Bold = 0

# ---------------------------------------------------------------------------
# 2) "Matthew Jovani ()" -> "Matthew Jovani (JOVM19108705)"
#    also removes the spell-check proofErr wrapper around "Jovani"
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("Matthew Jovani ()", $true, $false, $false, $false, $false, $true, 1, $false, "Matthew Jovani ()", 2)

$rng2b = $d.Content
$rng2b.Find.Execute("Matthew Jovani (")
$rng2b.Collapse(0)
$rng2b.InsertAfter("JOVM19108705")
$rng2b.Font.Bold = 1
$rng2b.Font.Bold = 0

# ---------------------------------------------------------------------------
# 3) Merge "Microsoft " / "Form" (proofErr-wrapped) / ". Le questionnaire..."
#    into a single run, removing both proofErr wrappers. The sentence before
#    it ("Visiteur comme employe ... corporelle.") must stay its own run, so
#    it is temporarily given different formatting to stop it from being
#    coalesced into the merged run, then restored.
# ---------------------------------------------------------------------------
$rngProtect = $d.Content
$rngProtect.Find.Execute("Visiteur comme employ" + [char]0x00E9 + " devront passer un contr" + [char]0x00F4 + "le de temp" + [char]0x00E9 + "rature corporelle.")
$rngProtect.Font.Bold = 1

$rng3 = $d.Content
$search3 = "Microsoft Form. Le questionnaire pr" + [char]0x00E9 + "sentera les questions suivantes" + [char]0x00A0 + ":"
$rng3.Find.Execute($search3, $true, $false, $false, $false, $false, $true, 1, $false, $search3, 2)

$rngUnprotect = $d.Content
$rngUnprotect.Find.Execute("Visiteur comme employ" + [char]0x00E9 + " devront passer un contr" + [char]0x00F4 + "le de temp" + [char]0x00E9 + "rature corporelle.")
$rngUnprotect.Font.Bold = 0

Write-Output "done"
